# "Add ui changes from meeting"
#
# 1) Refresh the cached "datetimeFigureOut" footer-date text (30.08.2019 ->
#    06.09.2019) everywhere it is cached: on the slide master and on every
#    slide layout's date placeholder.
# 2) Split the "Ranks Tab" button label on slide 1 into two runs so the
#    first word reads "Leaderboard" instead of "Ranks" ("Leaderboard Tab").

$p = $ppt.ActivePresentation

$oldDate = "30.08.2019"
$newDate = "06.09.2019"

function Update-DatePlaceholder {
    param($shapes)

    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master's date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's date placeholder.
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 1: "Ranks Tab" -> "Leaderboard Tab" (as two runs: "Leaderboard" + " Tab").
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes
for ($k = 1; $k -le $shapes.Count; $k++) {
    $sh = $shapes.Item($k)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "Ranks Tab") {
            $tr.Text = "Leaderboard"
            [void]$tr.InsertAfter(" Tab")
        }
    }
}
